# Apply the edits described by the diff:
#  1) Slide 2 ("Uge 43/44" schedule): three small wording tweaks inside the
#     Monday/Wednesday paragraph of the content placeholder.
#  2) Slide 5 ("Relationsdatabase"): give the title placeholder an explicit
#     position/size (xfrm) instead of inheriting it from the layout/master.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Text tweaks on slide 2
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
$tr = $sh2.TextFrame.TextRange

# -- Run A: "... se film omkring ER 	diagram evt. prøv det af"
#           -> "... evt. se film omkring ER diagram og prøv det af"
$oldA = "`" eller `"DB Browser`" eller MySQL server. 			Sætte en simpel database op, prøve simple SQL og evt. PHP kommandoer 			af, se film omkring ER 	diagram evt. prøv det af"
$newA = "” eller ”DB Browser” eller MySQL server. 			Sætte en simpel database op, prøve simple SQL og evt. PHP kommandoer 			af, evt. se film omkring ER diagram og prøv det af"
$idxA = $tr.Text.IndexOf($oldA)
if ($idxA -lt 0) { throw "Run A (DB Browser sentence) not found" }
$tr.Characters($idxA + 1, $oldA.Length).Text = $newA

# -- Run B: " Forsætte med at forstå SQL ..." -> " Fortsætte med at forstå SQL ..."
$oldB = " Forsætte med at forstå SQL og databaser samt prøve simple 			(udleverede) "
$newB = " Fortsætte med at forstå SQL og databaser samt prøve simple 			(udleverede) "
$idxB = $tr.Text.IndexOf($oldB)
if ($idxB -lt 0) { throw "Run B (Forsaette) not found" }
$tr.Characters($idxB + 1, $oldB.Length).Text = $newB

# -- Run C: "... Fortsæt med ER diagrammer." -> "... Fortsæt evt. med ER diagrammer."
$oldC = " programmer af imod databasen og få data frem og 			tilbage vha. simple SQL kommandoer. Fortsæt med ER diagrammer."
$newC = " programmer af imod databasen og få data frem og 			tilbage vha. simple SQL kommandoer. Fortsæt evt. med ER diagrammer."
$idxC = $tr.Text.IndexOf($oldC)
if ($idxC -lt 0) { throw "Run C (Fortsaet med ER) not found" }
$tr.Characters($idxC + 1, $oldC.Length).Text = $newC

# ---------------------------------------------------------------------
# 2) Position the title placeholder on slide 5
# ---------------------------------------------------------------------

# Shape.Left/Top/Width/Height are expressed in points (1 pt = 12700 EMU),
# but internally the value round-trips through a single-precision float
# before being re-quantised to EMU, so naively using targetEmu/12700.0
# can be off by 1 EMU. Nudge the point value in tiny steps until the
# float32 -> EMU pipeline reproduces the exact target.
function Get-PtForEmu {
    param([double]$TargetEmu)
    $base = $TargetEmu / 12700.0
    for ($i = 0; $i -le 1000; $i++) {
        foreach ($sign in 1, -1) {
            $cand = $base + ($sign * $i * 0.00001)
            $asFloat = [double]([float]$cand)
            $emu = [math]::Floor($asFloat * 12700)
            if ($emu -eq $TargetEmu) {
                return $cand
            }
        }
    }
    return $base
}

$s5 = $p.Slides.Item(5)
$titleShape = $s5.Shapes.Item(1)

$titleShape.Left = Get-PtForEmu 411480
$titleShape.Top = Get-PtForEmu 222885
$titleShape.Width = Get-PtForEmu 10515600
$titleShape.Height = Get-PtForEmu 1325563
